# chore: update Sheets via scheduled runner
# Refresh cached Universalis price snapshots (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ) and the dependent
# LevePriceNQ/HQ + LeveProfitNQ/HQ columns across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 179668.14
$ws.Range("J17").Value = 184631.14
$ws.Range("L17").Value = 553893.42
$ws.Range("N17").Value = -554229.42

$ws.Range("H28").Value = 568.875
$ws.Range("I28").Value = 599.1053000000001
$ws.Range("J28").Value = 454
$ws.Range("K28").Value = 599.1053000000001
$ws.Range("L28").Value = 454
$ws.Range("M28").Value = -114.1053000000001
$ws.Range("N28").Value = -1424

$ws.Range("H57").Value = 81666.664
$ws.Range("J57").Value = 81666.664
$ws.Range("L57").Value = 244999.992
$ws.Range("N57").Value = -245997.992

$ws.Range("H86").Value = 102567790
$ws.Range("I86").Value = 142860500
$ws.Range("J86").Value = 55559640
$ws.Range("K86").Value = 142860500
$ws.Range("L86").Value = 55559640
$ws.Range("M86").Value = -142859377
$ws.Range("N86").Value = -55561886

$ws.Range("H89").Value = 102567790
$ws.Range("I89").Value = 142860500
$ws.Range("J89").Value = 55559640
$ws.Range("K89").Value = 714302500
$ws.Range("L89").Value = 277798200
$ws.Range("M89").Value = -714296884
$ws.Range("N89").Value = -277809432

$ws.Range("H134").Value = 99975
$ws.Range("J134").Value = 99975
$ws.Range("L134").Value = 99975
$ws.Range("N134").Value = -110115

$ws.Range("H138").Value = 8611.434999999999
$ws.Range("J138").Value = 3173.7058
$ws.Range("L138").Value = 9521.117400000001
$ws.Range("N138").Value = -19801.1174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 404656.53
$ws.Range("I32").Value = 667593.5600000001
$ws.Range("J32").Value = 10251
$ws.Range("K32").Value = 667593.5600000001
$ws.Range("L32").Value = 10251
$ws.Range("M32").Value = -667306.5600000001
$ws.Range("N32").Value = -10825

$ws.Range("H45").Value = 57662.723
$ws.Range("I45").Value = 73315
$ws.Range("K45").Value = 73315
$ws.Range("M45").Value = -72938

$ws.Range("H46").Value = 7217.3335
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7217.3335
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 7217.3335
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -7855.3335

$ws.Range("H74").Value = 354797.56
$ws.Range("I74").Value = 2406.2888
$ws.Range("J74").Value = 1109921.8
$ws.Range("K74").Value = 2406.2888
$ws.Range("L74").Value = 1109921.8
$ws.Range("M74").Value = -1532.2888
$ws.Range("N74").Value = -1111669.8

$ws.Range("H77").Value = 354797.56
$ws.Range("I77").Value = 2406.2888
$ws.Range("J77").Value = 1109921.8
$ws.Range("K77").Value = 12031.444
$ws.Range("L77").Value = 5549609
$ws.Range("M77").Value = -7663.444
$ws.Range("N77").Value = -5558345

$ws.Range("H110").Value = 1747.375
$ws.Range("I110").Value = 1747.375
$ws.Range("K110").Value = 1747.375
$ws.Range("M110").Value = 297.625

$ws.Range("H140").Value = 82331.664
$ws.Range("J140").Value = 83997.5
$ws.Range("L140").Value = 83997.5
$ws.Range("N140").Value = -94357.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 31999.5
$ws.Range("I99").Value = 31999.5
$ws.Range("K99").Value = 31999.5
$ws.Range("M99").Value = -30501.5

$ws.Range("H107").Value = 8516.448
$ws.Range("I107").Value = 9167.385
$ws.Range("J107").Value = 2875
$ws.Range("K107").Value = 9167.385
$ws.Range("L107").Value = 2875
$ws.Range("M107").Value = -7247.385
$ws.Range("N107").Value = -6715

$ws.Range("H122").Value = 49999
$ws.Range("J122").Value = 49999
$ws.Range("L122").Value = 49999
$ws.Range("N122").Value = -59799

$ws.Range("H134").Value = 21952874
$ws.Range("I134").Value = 1335.5
$ws.Range("J134").Value = 60002210
$ws.Range("K134").Value = 4006.5
$ws.Range("L134").Value = 180006630
$ws.Range("M134").Value = -1471.5
$ws.Range("N134").Value = -180011700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 14287273
$ws.Range("I16").Value = 17858592
$ws.Range("J16").Value = 1995
$ws.Range("K16").Value = 17858592
$ws.Range("L16").Value = 1995
$ws.Range("M16").Value = -17858305
$ws.Range("N16").Value = -2569

$ws.Range("H31").Value = 2784.4517
$ws.Range("I31").Value = 3230.7334
$ws.Range("K31").Value = 3230.7334
$ws.Range("M31").Value = -2935.7334

$ws.Range("H34").Value = 2784.4517
$ws.Range("I34").Value = 3230.7334
$ws.Range("K34").Value = 3230.7334
$ws.Range("M34").Value = -3028.7334

$ws.Range("H55").Value = 3149.5
$ws.Range("I55").Value = 3149.5
$ws.Range("K55").Value = 3149.5
$ws.Range("M55").Value = -2834.5

$ws.Range("H107").Value = 2317.8948
$ws.Range("J107").Value = 2273.25
$ws.Range("L107").Value = 2273.25
$ws.Range("N107").Value = -6113.25

$ws.Range("H113").Value = 14287273
$ws.Range("I113").Value = 17858592
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 17858592
$ws.Range("L113").Value = 1995
$ws.Range("M113").Value = -17856422
$ws.Range("N113").Value = -6335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3498.8
$ws.Range("J68").Value = 3498
$ws.Range("L68").Value = 10494
$ws.Range("N68").Value = -12116

$ws.Range("H71").Value = 3498.8
$ws.Range("J71").Value = 3498
$ws.Range("L71").Value = 31482
$ws.Range("N71").Value = -39594

$ws.Range("H75").Value = 23814576
$ws.Range("J75").Value = 28577400
$ws.Range("L75").Value = 85732200
$ws.Range("N75").Value = -85734196

$ws.Range("H78").Value = 23814576
$ws.Range("J78").Value = 28577400
$ws.Range("L78").Value = 257196600
$ws.Range("N78").Value = -257206584

$ws.Range("H95").Value = 13330.2
$ws.Range("I95").Value = 3333
$ws.Range("K95").Value = 9999
$ws.Range("M95").Value = -7940

$ws.Range("H113").Value = 549.3333
$ws.Range("J113").Value = 549.3333
$ws.Range("L113").Value = 1647.9999
$ws.Range("N113").Value = -5987.9999

$ws.Range("H131").Value = 3345.0967
$ws.Range("J131").Value = 3378.1785
$ws.Range("L131").Value = 10134.5355
$ws.Range("N131").Value = -20214.5355

$ws.Range("H137").Value = 3554.4443
$ws.Range("J137").Value = 6500
$ws.Range("L137").Value = 19500
$ws.Range("N137").Value = -29700

$ws.Range("H138").Value = 4006.111
$ws.Range("I138").Value = 4154.067
$ws.Range("K138").Value = 12462.201
$ws.Range("M138").Value = -7322.201000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3898.875
$ws.Range("I70").Value = 3659.4
$ws.Range("J70").Value = 4298
$ws.Range("K70").Value = 3659.4
$ws.Range("L70").Value = 4298
$ws.Range("M70").Value = -3389.4
$ws.Range("N70").Value = -4838

$ws.Range("H73").Value = 3898.875
$ws.Range("I73").Value = 3659.4
$ws.Range("J73").Value = 4298
$ws.Range("K73").Value = 3659.4
$ws.Range("L73").Value = 4298
$ws.Range("M73").Value = -2723.4
$ws.Range("N73").Value = -6170

$ws.Range("H107").Value = 167832.5
$ws.Range("I107").Value = 999999
$ws.Range("K107").Value = 999999
$ws.Range("M107").Value = -998079

$ws.Range("H122").Value = 5607.154
$ws.Range("I122").Value = 5581.364
$ws.Range("K122").Value = 16744.092
$ws.Range("M122").Value = -14294.092

$ws.Range("H126").Value = 2659
$ws.Range("I126").Value = 2070.3333
$ws.Range("K126").Value = 6210.999899999999
$ws.Range("M126").Value = -3740.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2030.7273
$ws.Range("I61").Value = 1730.6842
$ws.Range("K61").Value = 1730.6842
$ws.Range("M61").Value = -1528.6842

$ws.Range("H68").Value = 3500
$ws.Range("J68").Value = 3500
$ws.Range("L68").Value = 3500
$ws.Range("N68").Value = -4998

$ws.Range("H71").Value = 3500
$ws.Range("J71").Value = 3500
$ws.Range("L71").Value = 17500
$ws.Range("N71").Value = -24988

$ws.Range("H113").Value = 2030.7273
$ws.Range("I113").Value = 1730.6842
$ws.Range("K113").Value = 1730.6842
$ws.Range("M113").Value = 439.3158000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5795511.5
$ws.Range("I81").Value = 6480584.5
$ws.Range("K81").Value = 12961169
$ws.Range("M81").Value = -12960108

$ws.Range("H84").Value = 5795511.5
$ws.Range("I84").Value = 6480584.5
$ws.Range("K84").Value = 64805845
$ws.Range("M84").Value = -64800541

$ws.Range("H107").Value = 1787698.4
$ws.Range("I107").Value = 1437.5
$ws.Range("K107").Value = 4312.5
$ws.Range("M107").Value = -2392.5

$ws.Range("H113").Value = 678.625
$ws.Range("I113").Value = 608.6
$ws.Range("J113").Value = 795.3333
$ws.Range("K113").Value = 1825.8
$ws.Range("L113").Value = 2385.9999
$ws.Range("M113").Value = 344.1999999999998
$ws.Range("N113").Value = -6725.9999

$ws.Range("H126").Value = 11364946
$ws.Range("I126").Value = 12501041
$ws.Range("K126").Value = 37503123
$ws.Range("M126").Value = -37500653
